{"js": "// Remove the trailing \"Ver no Jupiter ...\" / copyright footer block from the\n// end of the document, along with the blank paragraph that separates it from\n// the last requirement line (\"LOQ4053: ...\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\n\n// Locate the \"Ver no Jupiter...\" paragraph; the footer block consists of\n// that paragraph, the copyright paragraph right after it, and the blank\n// paragraph that immediately precedes them.\nlet verIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targets[0]) {\n    verIndex = i;\n    break;\n  }\n}\n\nif (verIndex !== -1) {\n  const toDelete = [];\n  // Blank paragraph right before \"Ver no Jupiter...\" (if present and empty).\n  if (verIndex - 1 >= 0 && items[verIndex - 1].text === \"\") {\n    toDelete.push(items[verIndex - 1]);\n  }\n  // The \"Ver no Jupiter...\" paragraph itself.\n  toDelete.push(items[verIndex]);\n  // The copyright paragraph right after it (if it matches).\n  if (verIndex + 1 < items.length && items[verIndex + 1].text === targets[1]) {\n    toDelete.push(items[verIndex + 1]);\n  }\n\n  for (const p of toDelete) {\n    p.delete();\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$verText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n# Find the \"Ver no Jupiter...\" paragraph (searching from the end, since it's\n# near the end of the document).\n$verIndex = -1\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $verText) {\n        $verIndex = $i\n        break\n    }\n}\n\nif ($verIndex -ge 1) {\n    $copyrightIndex = $verIndex + 1\n    $blankIndex = $verIndex - 1\n\n    $hasCopyright = $false\n    if ($copyrightIndex -le $d.Paragraphs.Count) {\n        $ct = $d.Paragraphs.Item($copyrightIndex).Range.Text.TrimEnd([char]13, [char]7)\n        if ($ct -eq $copyrightText) {\n            $hasCopyright = $true\n        }\n    }\n\n    $hasBlank = $false\n    if ($blankIndex -ge 1) {\n        $bt = $d.Paragraphs.Item($blankIndex).Range.Text.TrimEnd([char]13, [char]7)\n        if ($bt -eq \"\") {\n            $hasBlank = $true\n        }\n    }\n\n    # Delete from the highest index downward so earlier indices stay valid.\n    if ($hasCopyright) {\n        $d.Paragraphs.Item($copyrightIndex).Range.Delete()\n    }\n    $d.Paragraphs.Item($verIndex).Range.Delete()\n    if ($hasBlank) {\n        $d.Paragraphs.Item($blankIndex).Range.Delete()\n    }\n}\n"}
